# Update the two-digit multiplication answers to match the new output.
$d = $word.ActiveDocument

$pairs = @(
    @("83×53=4399", "60×87=5220"),
    @("93×90=8370", "60×86=5160"),
    @("56×24=1344", "39×75=2925"),
    @("18×62=1116", "84×85=7140"),
    @("97×57=5529", "96×49=4704"),
    @("55×74=4070", "15×65=975"),
    @("27×55=1485", "42×96=4032"),
    @("55×68=3740", "85×90=7650"),
    @("18×74=1332", "50×86=4300"),
    @("30×48=1440", "27×44=1188"),
    @("92×51=4692", "11×70=770"),
    @("31×66=2046", "65×91=5915"),
    @("85×63=5355", "70×41=2870"),
    @("27×40=1080", "17×14=238"),
    @("26×25=650",  "62×70=4340"),
    @("56×68=3808", "12×85=1020"),
    @("61×97=5917", "27×88=2376"),
    @("24×58=1392", "22×49=1078"),
    @("47×89=4183", "82×93=7626"),
    @("43×19=817",  "79×73=5767"),
    @("59×79=4661", "83×30=2490"),
    @("15×57=855",  "44×81=3564"),
    @("12×40=480",  "34×62=2108"),
    @("64×65=4160", "26×16=416"),
    @("27×28=756",  "39×77=3003")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
